$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.135.33'
$ws.Range('E2').Value = '  -1.92%  '
$ws.Range('D3').Value = '1.821.05'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.58%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.27'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.004'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.49%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4227'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.88%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3682'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07232'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8536'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.93'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.97%  '
$ws.Range('D12').Value = '1.825.18'
$ws.Range('E12').Value = '  -1.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.689'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.60%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.07084'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.69%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.289'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -3.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '89.17'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.61%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.006'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.67%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008828'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.74%  '
$ws.Range('E19').Value = '  -0.47%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.00'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -2.83%  '
$ws.Range('B21').Value = 'WrappedBTC'
$ws.Range('C21').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D21').Value = '27.216.65'
$ws.Range('E21').Value = '  -1.64%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.115'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.54%  '
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.85'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.68%  '
$ws.Range('B24').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C24').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D24').Value = '2.050.44'
$ws.Range('E24').Value = '  -0.97%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.984'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.08%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '152.41'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.03%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.190'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.48%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.37'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.05%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.213'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -3.20%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.90'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -3.54%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08834'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.96%  '
$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.185'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.57%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7465'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -4.08%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.970'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.13%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.434'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.70%  '
$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.005'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.46%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.112'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.11%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01962'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.38%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05234'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.74%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.294'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.86%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.877'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.31%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1692'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.93%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5020'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.68%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.636'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -3.13%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.58'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.42%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '106.21'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.33%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4729'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.004'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.47%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06390'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.73%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.658'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.35%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.861'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.34%  '
